$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Test 7: Unittest isScaleneTest()) där en sida större än summan av de två andra ---
# Initial run (before Buggfix 3) - bug is still present, test result wrongly shows Pass
$ws.Cells.Item(28, 1).Value = 7
$ws.Cells.Item(28, 2).Value = "Unittest isScaleneTest()) där en sida större än summan av de två andra"
$ws.Cells.Item(28, 3).Value = "1.2, 20, 1.3"
$ws.Cells.Item(28, 4).Value = "Fail"
$ws.Cells.Item(28, 5).Value = "Pass"

# Description of the bug fix applied
$ws.Cells.Item(30, 2).Value = "Buggfix 3 - kontroll av inparametrar i konstruktorn (a + b <= c) || (a + c <= b) || (b + c <= a)"

# Re-run of test 7 after Buggfix 3 - now correctly shows Fail
$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "Unittest isScaleneTest()) där en sida större än summan av de två andra"
$ws.Cells.Item(32, 3).Value = "1.2, 20, 1.3"
$ws.Cells.Item(32, 4).Value = "Fail"
$ws.Cells.Item(32, 5).Value = "Fail"

# Scroll / selection state after the edit, matching the author's view when saving
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("E32").Select()
